$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "some text from Excel"
$ws.Range("B1").Value = "some value from Excel"
$ws.Range("A2").Value = "some text from Excel"
$ws.Range("B2").Value = "some value from Excel"
